$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the helper/example text under "Nama Posisi *" (column B, row 2):
# old text talked about allowing multiple comma-separated positions,
# new text restricts initialization to a single main position.
$ws.Range("B2").Value = "(Wajib sama dengan data dari database, hanya boleh mengisi 1 Posisi utama sebagai inisialisasi )"

# Widen column B (best-fit) so the longer instruction text fits.
$ws.Columns.Item(2).ColumnWidth = 95.25

# Move the active selection/view back to the start of the sheet (B5),
# instead of being scrolled over to column AB with AI4 selected.
$ws.Range("B5").Select()
